$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for the new worker (pushes the signature block down) ---
$ws.Rows.Item(17).Insert()

# Copy formatting (borders/fonts/number formats) from the existing data row (16)
# into the freshly inserted row (17) so it matches the table styling.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Populate the new worker row ---
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1002440552"
$ws.Range("D17").Value = "JONATHAN ALFREDO ORTEGA TAPIA"
$ws.Range("E17").Value = "2405"
$ws.Range("F17").Value = 12133
$ws.Range("G17").Value = 1300000

# --- Update the summary totals at the top of the statement ---
$ws.Range("E11").Value = 70800
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2

# --- Widen column D so the longer new name fits ---
$ws.Columns.Item(4).ColumnWidth = 33.453125
